$wb = $excel.ActiveWorkbook

# Sheet ALC, row 6
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 1759.3077
$ws.Cells.Item(6, 9).Value = 1889.3334
$ws.Cells.Item(6, 11).Value = 5668.0002
$ws.Cells.Item(6, 13).Value = -5556.0002

# Sheet ALC, row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 1399.6
$ws.Cells.Item(19, 9).Value = 1171.1428
$ws.Cells.Item(19, 10).Value = 1932.6666
$ws.Cells.Item(19, 11).Value = 1171.1428
$ws.Cells.Item(19, 12).Value = 1932.6666
$ws.Cells.Item(19, 13).Value = -996.1428000000001
$ws.Cells.Item(19, 14).Value = -2282.6666

# Sheet ALC, row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 295.2857
$ws.Cells.Item(41, 9).Value = 233.6
$ws.Cells.Item(41, 11).Value = 233.6
$ws.Cells.Item(41, 13).Value = 206.4

# Sheet ALC, row 48
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(48, 8).Value = 9475.333000000001
$ws.Cells.Item(48, 9).Value = 0
$ws.Cells.Item(48, 10).Value = 9475.333000000001
$ws.Cells.Item(48, 11).Value = 0
$ws.Cells.Item(48, 12).Value = 28425.999
$ws.Cells.Item(48, 14).Value = -29009.999
$ws.Cells.Item(48, 13).ClearContents()

# Sheet ALC, row 56
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(56, 8).Value = 9475.333000000001
$ws.Cells.Item(56, 9).Value = 0
$ws.Cells.Item(56, 10).Value = 9475.333000000001
$ws.Cells.Item(56, 11).Value = 0
$ws.Cells.Item(56, 12).Value = 28425.999
$ws.Cells.Item(56, 14).Value = -29493.999
$ws.Cells.Item(56, 13).ClearContents()

# Sheet ALC, row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 3910.1667
$ws.Cells.Item(62, 9).Value = 3547
$ws.Cells.Item(62, 11).Value = 3547
$ws.Cells.Item(62, 13).Value = -2923

# Sheet ALC, row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 3910.1667
$ws.Cells.Item(65, 9).Value = 3547
$ws.Cells.Item(65, 11).Value = 17735
$ws.Cells.Item(65, 13).Value = -14615

# Sheet ALC, row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 2135.818
$ws.Cells.Item(100, 9).Value = 2161.875
$ws.Cells.Item(100, 11).Value = 2161.875
$ws.Cells.Item(100, 13).Value = -1620.875

# Sheet ALC, row 115
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(115, 8).Value = 377
$ws.Cells.Item(115, 9).Value = 377
$ws.Cells.Item(115, 11).Value = 1131
$ws.Cells.Item(115, 13).Value = 436

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1663.3334
$ws.Cells.Item(137, 10).Value = 4005
$ws.Cells.Item(137, 12).Value = 12015
$ws.Cells.Item(137, 14).Value = -17115

# Sheet ARM, row 26
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(26, 8).Value = 15335.333
$ws.Cells.Item(26, 9).Value = 10003
$ws.Cells.Item(26, 10).Value = 26000
$ws.Cells.Item(26, 11).Value = 10003
$ws.Cells.Item(26, 12).Value = 26000
$ws.Cells.Item(26, 13).Value = -9673
$ws.Cells.Item(26, 14).Value = -26660

# Sheet ARM, row 29
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(29, 8).Value = 49999
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 49999
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 12).Value = 49999
$ws.Cells.Item(29, 14).Value = -50615
$ws.Cells.Item(29, 13).ClearContents()

# Sheet ARM, row 30
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 14).ClearContents()
$ws.Cells.Item(30, 13).ClearContents()

# Sheet ARM, row 49
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(49, 8).Value = 52000
$ws.Cells.Item(49, 10).Value = 52000
$ws.Cells.Item(49, 12).Value = 52000
$ws.Cells.Item(49, 14).Value = -52520

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1661.1538
$ws.Cells.Item(132, 9).Value = 1609.091
$ws.Cells.Item(132, 11).Value = 4827.272999999999
$ws.Cells.Item(132, 13).Value = -2297.272999999999

# Sheet BSM, row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 4374.75
$ws.Cells.Item(20, 9).Value = 4000
$ws.Cells.Item(20, 11).Value = 4000
$ws.Cells.Item(20, 13).Value = -3753

# Sheet BSM, row 29
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(29, 8).Value = 999.5
$ws.Cells.Item(29, 9).Value = 999.5
$ws.Cells.Item(29, 11).Value = 999.5
$ws.Cells.Item(29, 13).Value = -710.5

# Sheet BSM, row 30
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 14).ClearContents()

# Sheet BSM, row 44
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(44, 8).Value = 20000
$ws.Cells.Item(44, 10).Value = 20000
$ws.Cells.Item(44, 12).Value = 20000
$ws.Cells.Item(44, 14).Value = -20994

# Sheet BSM, row 64
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 2006.3334
$ws.Cells.Item(64, 9).Value = 2506
$ws.Cells.Item(64, 10).Value = 1007
$ws.Cells.Item(64, 11).Value = 2506
$ws.Cells.Item(64, 12).Value = 1007
$ws.Cells.Item(64, 13).Value = -2281
$ws.Cells.Item(64, 14).Value = -1457

# Sheet BSM, row 67
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(67, 8).Value = 2006.3334
$ws.Cells.Item(67, 9).Value = 2506
$ws.Cells.Item(67, 10).Value = 1007
$ws.Cells.Item(67, 11).Value = 2506
$ws.Cells.Item(67, 12).Value = 1007
$ws.Cells.Item(67, 13).Value = -1726
$ws.Cells.Item(67, 14).Value = -2567

# Sheet BSM, row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 5449.75
$ws.Cells.Item(99, 9).Value = 5933
$ws.Cells.Item(99, 11).Value = 5933
$ws.Cells.Item(99, 13).Value = -4435

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 6372.375
$ws.Cells.Item(134, 9).Value = 6473.2383
$ws.Cells.Item(134, 10).Value = 5666.3335
$ws.Cells.Item(134, 11).Value = 19419.7149
$ws.Cells.Item(134, 12).Value = 16999.0005
$ws.Cells.Item(134, 13).Value = -16884.7149
$ws.Cells.Item(134, 14).Value = -22069.0005

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 1950
$ws.Cells.Item(132, 9).Value = 1950
$ws.Cells.Item(132, 11).Value = 5850
$ws.Cells.Item(132, 13).Value = -3320

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 2710.3333
$ws.Cells.Item(134, 9).Value = 2710.3333
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 8130.999899999999
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 13).Value = -5595.999899999999
$ws.Cells.Item(134, 14).ClearContents()

# Sheet CUL, row 6
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(6, 8).Value = 5668.1665
$ws.Cells.Item(6, 9).Value = 501.75
$ws.Cells.Item(6, 11).Value = 1505.25
$ws.Cells.Item(6, 13).Value = -1392.25

# Sheet GSM, row 63
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(63, 8).Value = 30551.5
$ws.Cells.Item(63, 9).Value = 11103
$ws.Cells.Item(63, 10).Value = 50000
$ws.Cells.Item(63, 11).Value = 11103
$ws.Cells.Item(63, 12).Value = 50000
$ws.Cells.Item(63, 13).Value = -10417
$ws.Cells.Item(63, 14).Value = -51372

# Sheet GSM, row 66
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(66, 8).Value = 30551.5
$ws.Cells.Item(66, 9).Value = 11103
$ws.Cells.Item(66, 10).Value = 50000
$ws.Cells.Item(66, 11).Value = 33309
$ws.Cells.Item(66, 12).Value = 150000
$ws.Cells.Item(66, 13).Value = -29877
$ws.Cells.Item(66, 14).Value = -156864

# Sheet GSM, row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 819.8570999999999
$ws.Cells.Item(97, 9).Value = 819.8570999999999
$ws.Cells.Item(97, 11).Value = 819.8570999999999
$ws.Cells.Item(97, 13).Value = -323.8570999999999

# Sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2599
$ws.Cells.Item(102, 9).Value = 2599
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 2599
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = -977
$ws.Cells.Item(102, 14).ClearContents()

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 2162.3076
$ws.Cells.Item(122, 9).Value = 2092.5
$ws.Cells.Item(122, 10).Value = 3000
$ws.Cells.Item(122, 11).Value = 6277.5
$ws.Cells.Item(122, 12).Value = 9000
$ws.Cells.Item(122, 13).Value = -3827.5
$ws.Cells.Item(122, 14).Value = -13900

# Sheet GSM, row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 16699.75
$ws.Cells.Item(126, 9).Value = 16699.75
$ws.Cells.Item(126, 11).Value = 50099.25
$ws.Cells.Item(126, 13).Value = -47629.25

# Sheet LTW, row 74
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(74, 8).Value = 41248.5
$ws.Cells.Item(74, 10).Value = 41248.5
$ws.Cells.Item(74, 12).Value = 41248.5
$ws.Cells.Item(74, 14).Value = -43244.5

# Sheet LTW, row 77
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(77, 8).Value = 41248.5
$ws.Cells.Item(77, 10).Value = 41248.5
$ws.Cells.Item(77, 12).Value = 123745.5
$ws.Cells.Item(77, 14).Value = -133729.5

# Sheet WVR, row 4
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 1550.125
$ws.Cells.Item(4, 9).Value = 366.66666
$ws.Cells.Item(4, 11).Value = 366.66666
$ws.Cells.Item(4, 13).Value = -253.66666

# Sheet WVR, row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 252.16667
$ws.Cells.Item(113, 9).Value = 201
$ws.Cells.Item(113, 10).Value = 277.75
$ws.Cells.Item(113, 11).Value = 603
$ws.Cells.Item(113, 12).Value = 833.25
$ws.Cells.Item(113, 13).Value = 1567
$ws.Cells.Item(113, 14).Value = -5173.25
